# PXK.xlsx template edit: swap placeholder tokens to the new ImportExcel-era
# variable names (outDeptName/outDeptAddress/e.index/e.name/e.code/e.unit/e.realQuantity).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

# --- Header block (rows 9-13): sender/department info now reuses the same
#     "outDept*" fields for both the "Người nhận hàng" and "Xuất tại kho" rows.
$ws.Range("C9").Value  = "`${outDeptName}"
$ws.Range("C10").Value = "`${outDeptAddress}"
$ws.Range("C12").Value = "`${outDeptName}"
$ws.Range("C13").Value = "`${outDeptAddress}"

# --- Table header row 17 (data row template for each line item "e")
$ws.Range("A17").Value = "`${e.index} "
$ws.Range("B17").Value = "`${e.name}"
$ws.Range("I17").Value = "`${e.code}"
$ws.Range("J17").Value = "`${e.unit}"
$ws.Range("L17").Value = "`${e.realQuantity}"

$ws.Range("A14").Select()
